$wb = $excel.ActiveWorkbook

# Add the new "estimates" worksheet as the last sheet (it was inserted
# right after Table 2 / "redo figures" commit adds a 3rd results tab).
$ws = $wb.Worksheets.Add()
$ws.Name = "estimates"
$ws.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Re-fetch by name: the handle returned by Add()/Move() can go stale
# once the sheet collection is restructured.
$ws = $wb.Worksheets("estimates")

# Left block: "JOINED" estimates.
$ws.Range("B3").Value = "JOINED"
$ws.Range("C4").Value = "b=67.89"
$ws.Range("C5").Value = "p=0.025"
$ws.Range("B6").Value = "JOINED + ALWAYS"
$ws.Range("C7").Value = "b=24"
$ws.Range("C8").Value = "p=.19"

# Right block: "LEFT" estimates.
$ws.Range("E3").Value = "LEFT"
$ws.Range("F4").Value = "b=16.28"
$ws.Range("F5").Value = "p=0.008"
$ws.Range("E6").Value = "LEFT +ALWAYS"
$ws.Range("F7").Value = "b=6.83"
$ws.Range("F8").Value = "p=.209"

$null = $ws.Range("E9").Select()
$null = $wb.Worksheets("estimates").Activate()
